# Chatflix QA Defect Report — the QA reviewer re-tested the remaining
# "Open" (visible, unfiltered) defects — rows 8-12 — and marked their
# Status column (E) as "Fixed", then left the selection on the
# Description cell of the first of those rows (D8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E8:E12").Value = "Fixed"

$ws.Range("D8").Select()
